$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule table header row 7: the "Amount" condition cell (E7) had its text
# re-cased from "Amount > $param" to "amount > $param", and the cell's
# font was switched to Arial.
$cell = $ws.Range("E7")
$cell.Value = "amount > `$param"
$cell.Font.Name = "Arial"
